# Vega Monumental Concepción - Acelga: weekly refresh of the "Fecha" (and
# associated price-record) block in rows 80..205.
#
# The data in that range is organised as "records": a market-day
# observation is either a single row (an odd one-off "Primera"-only
# quality reading) or a pair of rows (a "Primera" row immediately
# followed by a "Segunda" row sharing the same Fecha). The weekly refresh
# pushes every existing record one slot later, inserts a brand-new record
# at the front (a normal Primera/Segunda pair dated serial 44579), and
# appends a duplicate of the record that used to be last (rows 204/205,
# serial 44512) as the new final record - growing the sheet from R205 to
# R207.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 80
$lastRow = 205
$numCols = 18   # columns A..R

# ---- 1. Read the existing rows 80..205 into memory -----------------
$data = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $row = @()
    for ($c = 1; $c -le $numCols; $c++) {
        $row += ,($ws.Cells.Item($r, $c).Value2)
    }
    $data += ,$row
}

# ---- 2. Group the rows into "records" -------------------------------
# A record is 2 rows if row i and row i+1 share the same Fecha (col 4)
# and row i is quality "Primera" (col 9) while row i+1 is "Segunda";
# otherwise the record is just the single row i.
$records = @()
$i = 0
while ($i -lt $data.Count) {
    $row0 = $data[$i]
    $isPair = $false
    if ($i + 1 -lt $data.Count) {
        $row1 = $data[$i + 1]
        if (($row1[3] -eq $row0[3]) -and ($row0[8] -eq "Primera") -and ($row1[8] -eq "Segunda")) {
            $isPair = $true
        }
    }
    $rec = @()
    if ($isPair) {
        $rec += ,$data[$i]
        $rec += ,$data[$i + 1]
        $records += ,$rec
        $i += 2
    } else {
        $rec += ,$data[$i]
        $records += ,$rec
        $i += 1
    }
}

# ---- 3. Build the brand-new record (standard Primera/Segunda pair) --
# Columns A,B,C,E,F,G,H,R are constant throughout the block, so clone
# them from the existing first record; only Fecha/I/J/K/L/M/N/O/P/Q
# differ between a "standard" Primera row and a "standard" Segunda row.
$template = $records[0]
$newPrimera = $template[0]
$newSegunda = $template[1]
$newPrimera[3] = 44579
$newSegunda[3] = 44579
$newRecord = @()
$newRecord += ,$newPrimera
$newRecord += ,$newSegunda

# ---- 4. Re-sequence: new record first, then all-but-the-last of the
#         old records (shifting them one slot later), then a duplicate
#         of the old last record appended at the very end. -----------
$newRecords = @()
$newRecords += ,$newRecord
for ($k = 0; $k -lt $records.Count - 1; $k++) {
    $newRecords += ,$records[$k]
}
$newRecords += ,$records[$records.Count - 1]

# ---- 5. Flatten back out into physical rows starting at row 80 ------
$outRows = @()
foreach ($rec in $newRecords) {
    foreach ($r in $rec) {
        $outRows += ,$r
    }
}

$r = $firstRow
foreach ($row in $outRows) {
    for ($c = 1; $c -le $numCols; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $row[$c - 1]
    }
    # Preserve the date number format on column D (style index "2" in
    # the original file) for every row, including the two brand-new
    # ones past the old end (206/207).
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $r++
}
